$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 270
$ws1.Range("F4").Value = 1100
$ws1.Range("F5").Value = 2608
$ws1.Range("F7").Value = 682
$ws1.Range("F8").Value = 57
$ws1.Range("F9").Value = 242
$ws1.Range("F10").Value = 186
$ws1.Range("F11").Value = 687
$ws1.Range("F12").Value = 92
$ws1.Range("F13").Value = 115
$ws1.Range("F14").Value = 1504
$ws1.Range("F15").Value = 304

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("C4").Value = "广州·东方PartyNight×东方同人only-游剧天P2"
$ws2.Range("F6").Value = 14
$ws2.Range("F10").Value = 21
$ws2.Range("F12").Value = 43

# --- Sheet "本地生活" (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6347
$ws3.Range("F4").Value = 2014
$ws3.Range("F5").Value = 247

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6347
$ws4.Range("F4").Value = 2014
$ws4.Range("F5").Value = 247
$ws4.Range("C9").Value = "广州·东方PartyNight×东方同人only-游剧天P2"
$ws4.Range("F11").Value = 270
$ws4.Range("F12").Value = 1100
$ws4.Range("F13").Value = 14
$ws4.Range("F16").Value = 2608
$ws4.Range("F19").Value = 21
$ws4.Range("F21").Value = 43
$ws4.Range("F22").Value = 682
$ws4.Range("F23").Value = 57
$ws4.Range("F24").Value = 243
$ws4.Range("F26").Value = 186
$ws4.Range("F27").Value = 687
$ws4.Range("F28").Value = 92
$ws4.Range("F29").Value = 115
$ws4.Range("F31").Value = 1504
$ws4.Range("F32").Value = 304
